$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Notes!" label under the SHIFT block header (becomes shared string #36)
$ws.Range("A5").Value = "Notes!"

# Insert a new blank row above the old row 6 (shifts rows 6-16 down to 7-17,
# picking up formatting from the row above as Excel normally does)
$null = $ws.Rows.Item(6).Insert()

# Update the active selection left in the sheet view
$null = $ws.Range("B4").Select()
